$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "jobTitle"
$ws.Range("B1").Select() | Out-Null
